# Applies the commit "working with SCD30, PMS7003 and display":
#  1) In the `MySerial.begin(...)` paragraph, drop the spurious
#     grammar-check markers around the "(" run, merge the "(" run and
#     the "9600, SERIAL_8N1, 16, 17);" run into a single run, append a
#     new " // RX,TX" run (same Courier New / sz20 formatting), and move
#     the "_GoBack" bookmark onto the end of this paragraph.
#  2) Remove the (now relocated) "_GoBack" bookmark from the end of the
#     "...esp32-pinout-reference-gpios/):" paragraph.

$d = $word.ActiveDocument

function Find-ParagraphIndex($doc, [string]$needle) {
    for ($i = 1; $i -le $doc.Paragraphs.Count; $i++) {
        if ($doc.Paragraphs.Item($i).Range.Text -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# --- Paragraph 1: "    MySerial.begin(9600, SERIAL_8N1, 16, 17);" ---
$idx1 = Find-ParagraphIndex $d "SERIAL_8N1"
$p1 = $d.Paragraphs.Item($idx1)
$r1 = $p1.Range

$frag1 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00F909E2" w:rsidRPr="00F909E2" w:rsidRDefault="00F909E2" w:rsidP="00F909E2"><w:pPr><w:tabs><w:tab w:val="left" w:pos="916"/><w:tab w:val="left" w:pos="1832"/><w:tab w:val="left" w:pos="2748"/><w:tab w:val="left" w:pos="3664"/><w:tab w:val="left" w:pos="4580"/><w:tab w:val="left" w:pos="5496"/><w:tab w:val="left" w:pos="6412"/><w:tab w:val="left" w:pos="7328"/><w:tab w:val="left" w:pos="8244"/><w:tab w:val="left" w:pos="9160"/><w:tab w:val="left" w:pos="10076"/><w:tab w:val="left" w:pos="10992"/><w:tab w:val="left" w:pos="11908"/><w:tab w:val="left" w:pos="12824"/><w:tab w:val="left" w:pos="13740"/><w:tab w:val="left" w:pos="14656"/></w:tabs><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="de-DE"/></w:rPr></w:pPr><w:r w:rsidRPr="00F909E2"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve">    </w:t></w:r><w:proofErr w:type="spellStart"/><w:r w:rsidRPr="00F909E2"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="de-DE"/></w:rPr><w:t>MySerial.begin</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r w:rsidRPr="00F909E2"><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="de-DE"/></w:rPr><w:t>(9600, SERIAL_8N1, 16, 17);</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Courier New" w:eastAsia="Times New Roman" w:hAnsi="Courier New" w:cs="Courier New"/><w:sz w:val="20"/><w:szCs w:val="20"/><w:lang w:val="en-US" w:eastAsia="de-DE"/></w:rPr><w:t xml:space="preserve"> // RX,TX</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r1.InsertXML($frag1)

# --- Paragraph 2: "...esp32-pinout-reference-gpios/):" (drop the old bookmark) ---
$idx2 = Find-ParagraphIndex $d "esp32-pinout-reference-gpios"
$p2 = $d.Paragraphs.Item($idx2)
$r2 = $p2.Range

$frag2 = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="00F909E2" w:rsidRDefault="00F909E2" w:rsidP="00A917B5"><w:r><w:t xml:space="preserve">Ich habe gesehen dass UART0 auf Pins 1 (TX) und 3 (RX), sowie UART2 auf Pins 17 (TX)/16 (RX) gemappt wurden. Pins 34 bis 39 sind &#252;brigens ausschlie&#223;lich Eing&#228;nge. Komplett freie Pins f&#252;r I/O sind 4,5, und </w:t></w:r><w:r w:rsidR="005F68E3"><w:t xml:space="preserve">13-33 (von </w:t></w:r><w:r w:rsidR="005F68E3" w:rsidRPr="005F68E3"><w:t>https://randomnerdtutorials.com/esp32-pinout-reference-gpios/</w:t></w:r><w:r w:rsidR="005F68E3"><w:t>):</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

$r2.InsertXML($frag2)

Write-Output "applied edits to paragraphs $idx1 and $idx2"
